$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need NumberFormat "@"
# forced first, otherwise Excel auto-converts them to numeric cells and
# the literal text formatting (e.g. "1.000", "0.07689") would be lost.
$textFormatCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D17",
    "D18",
    "D19",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D40",
    "D41",
    "D42",
    "D43",
    "D45",
    "D46",
    "D48",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped by the GitHub Actions cron job.
$ws.Range("D2").Value2 = '26.425.99'
$ws.Range("E2").Value2 = '  +1.02%  '
$ws.Range("D3").Value2 = '1.844.11'
$ws.Range("E3").Value2 = '  -0.12%  '
$ws.Range("D4").Value2 = '1.000'
$ws.Range("E4").Value2 = '  +0.08%  '
$ws.Range("D5").Value2 = '259.80'
$ws.Range("E5").Value2 = '  -7.12%  '
$ws.Range("D6").Value2 = '1.000'
$ws.Range("E6").Value2 = '  +0.05%  '
$ws.Range("D7").Value2 = '0.5111'
$ws.Range("E7").Value2 = '  +0.02%  '
$ws.Range("D8").Value2 = '0.3229'
$ws.Range("E8").Value2 = '  -7.88%  '
$ws.Range("D9").Value2 = '0.06733'
$ws.Range("E9").Value2 = '  -1.50%  '
$ws.Range("D10").Value2 = '19.17'
$ws.Range("E10").Value2 = '  -3.99%  '
$ws.Range("D11").Value2 = '0.7712'
$ws.Range("E11").Value2 = '  -4.45%  '
$ws.Range("B12").Value2 = 'TRON'
$ws.Range("C12").Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value2 = '0.07689'
$ws.Range("E12").Value2 = '  -1.17%  '
$ws.Range("B13").Value2 = 'WrappedEther'
$ws.Range("C13").Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value2 = '1.877.28'
$ws.Range("E13").Value2 = '  +1.65%  '
$ws.Range("D14").Value2 = '88.10'
$ws.Range("E14").Value2 = '  -0.53%  '
$ws.Range("D15").Value2 = '5.025'
$ws.Range("E15").Value2 = '  -1.43%  '
$ws.Range("E16").Value2 = '  +0.20%  '
$ws.Range("D17").Value2 = '14.08'
$ws.Range("E17").Value2 = '  -0.84%  '
$ws.Range("D18").Value2 = '1.000'
$ws.Range("E18").Value2 = '  +0.07%  '
$ws.Range("D19").Value2 = '0.000007881'
$ws.Range("E19").Value2 = '  -2.44%  '
$ws.Range("D20").Value2 = '26.460.85'
$ws.Range("D21").Value2 = '2.119.16'
$ws.Range("E21").Value2 = '  +1.77%  '
$ws.Range("D22").Value2 = '4.573'
$ws.Range("E22").Value2 = '  -4.23%  '
$ws.Range("D23").Value2 = '9.543'
$ws.Range("E23").Value2 = '  -5.30%  '
$ws.Range("D24").Value2 = '5.987'
$ws.Range("E24").Value2 = '  -3.58%  '
$ws.Range("D25").Value2 = '2.343'
$ws.Range("E25").Value2 = '  -1.67%  '
$ws.Range("D26").Value2 = '145.34'
$ws.Range("E26").Value2 = '  +0.64%  '
$ws.Range("D27").Value2 = '1.651'
$ws.Range("E27").Value2 = '  -0.54%  '
$ws.Range("D28").Value2 = '16.91'
$ws.Range("E28").Value2 = '  -1.87%  '
$ws.Range("D29").Value2 = '110.74'
$ws.Range("E29").Value2 = '  +0.53%  '
$ws.Range("D30").Value2 = '4.211'
$ws.Range("E30").Value2 = '  -3.72%  '
$ws.Range("D31").Value2 = '4.178'
$ws.Range("E31").Value2 = '  -3.09%  '
$ws.Range("D32").Value2 = '0.08705'
$ws.Range("D33").Value2 = '0.04818'
$ws.Range("E33").Value2 = '  -1.89%  '
$ws.Range("E34").Value2 = '  -3.33%  '
$ws.Range("D35").Value2 = '2.862'
$ws.Range("E35").Value2 = '  +0.76%  '
$ws.Range("D36").Value2 = '0.6870'
$ws.Range("E36").Value2 = '  -6.93%  '
$ws.Range("D37").Value2 = '3.079'
$ws.Range("E37").Value2 = '  -4.89%  '
$ws.Range("E38").Value2 = '  -2.66%  '
$ws.Range("E39").Value2 = '  -6.77%  '
$ws.Range("D40").Value2 = '0.4926'
$ws.Range("E40").Value2 = '  -4.84%  '
$ws.Range("D41").Value2 = '113.16'
$ws.Range("E41").Value2 = '  -2.65%  '
$ws.Range("D42").Value2 = '0.9058'
$ws.Range("E42").Value2 = '  -6.11%  '
$ws.Range("D43").Value2 = '6.092'
$ws.Range("E43").Value2 = '  -2.74%  '
$ws.Range("D45").Value2 = '7.777'
$ws.Range("E45").Value2 = '  -2.97%  '
$ws.Range("D46").Value2 = '0.4263'
$ws.Range("E46").Value2 = '  -5.99%  '
$ws.Range("E47").Value2 = '  -5.61%  '
$ws.Range("D48").Value2 = '9.132'
$ws.Range("E48").Value2 = '  -2.46%  '
$ws.Range("D49").Value2 = '0.05903'
$ws.Range("E49").Value2 = '  -0.34%  '
$ws.Range("D50").Value2 = '35.18'
$ws.Range("D51").Value2 = '1.434'
$ws.Range("E51").Value2 = '  -4.51%  '
